$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 4 (Leve Item ID 5470)
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -728
# row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 99116.17999999999
$ws.Range("I51").Value = 207512.6
$ws.Range("J51").Value = 8785.833000000001
$ws.Range("K51").Value = 207512.6
$ws.Range("L51").Value = 8785.833000000001
$ws.Range("M51").Value = -207028.6
$ws.Range("N51").Value = -9753.833000000001
# row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 1466.6072
$ws.Range("I70").Value = 1798.3334
$ws.Range("J70").Value = 1309.4736
$ws.Range("K70").Value = 5395.0002
$ws.Range("L70").Value = 3928.4208
$ws.Range("M70").Value = -5125.0002
$ws.Range("N70").Value = -4468.4208
# row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 1466.6072
$ws.Range("I73").Value = 1798.3334
$ws.Range("J73").Value = 1309.4736
$ws.Range("K73").Value = 5395.0002
$ws.Range("L73").Value = 3928.4208
$ws.Range("M73").Value = -4459.0002
$ws.Range("N73").Value = -5800.4208
# row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 4884.5386
$ws.Range("I74").Value = 3499
$ws.Range("K74").Value = 3499
$ws.Range("M74").Value = -2563
# row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 4884.5386
$ws.Range("I77").Value = 3499
$ws.Range("K77").Value = 17495
$ws.Range("M77").Value = -12815
# row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 125064264
$ws.Range("I98").Value = 142930110
$ws.Range("K98").Value = 142930110
$ws.Range("M98").Value = -142928612
# row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 125064264
$ws.Range("I122").Value = 142930110
$ws.Range("K122").Value = 428790330
$ws.Range("M122").Value = -428787880

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 4 (Leve Item ID 5071)
$ws.Range("H4").Value = 1646.4
$ws.Range("J4").Value = 2998
$ws.Range("L4").Value = 2998
$ws.Range("N4").Value = -3230
# row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2028.5454
$ws.Range("I122").Value = 1101.8125
$ws.Range("K122").Value = 3305.4375
$ws.Range("M122").Value = -855.4375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 138
$ws.Range("I22").Value = 138
$ws.Range("K22").Value = 138
$ws.Range("M22").Value = 35
# row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 10876.5
$ws.Range("I99").Value = 13141.111
$ws.Range("J99").Value = 4082.6667
$ws.Range("K99").Value = 13141.111
$ws.Range("L99").Value = 4082.6667
$ws.Range("M99").Value = -11643.111
$ws.Range("N99").Value = -7078.6667
# row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 71752.734
$ws.Range("I134").Value = 1355.25
$ws.Range("K134").Value = 4065.75
$ws.Range("M134").Value = -1530.75
# row 141 (Leve Item ID 43278)
$ws.Range("H141").Value = 48777.5
$ws.Range("J141").Value = 48777.5
$ws.Range("L141").Value = 48777.5
$ws.Range("N141").Value = -59137.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 369.125
$ws.Range("I22").Value = 369.125
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 369.125
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -19.125
$ws.Range("N22").ClearContents()
# row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 4198.6
$ws.Range("I86").Value = 4198.6
$ws.Range("K86").Value = 4198.6
$ws.Range("M86").Value = -3075.6
# row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 4198.6
$ws.Range("I89").Value = 4198.6
$ws.Range("K89").Value = 20993
$ws.Range("M89").Value = -15377

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 73.84375
$ws.Range("I2").Value = 61.86207
$ws.Range("K2").Value = 371.17242
$ws.Range("M2").Value = -258.17242
# row 22 (Leve Item ID 4697)
$ws.Range("H22").Value = 2144.3076
$ws.Range("I22").Value = 404.33334
$ws.Range("J22").Value = 3635.7144
$ws.Range("K22").Value = 1213.00002
$ws.Range("L22").Value = 10907.1432
$ws.Range("M22").Value = -1044.00002
$ws.Range("N22").Value = -11245.1432
# row 27 (Leve Item ID 4697)
$ws.Range("H27").Value = 2144.3076
$ws.Range("I27").Value = 404.33334
$ws.Range("J27").Value = 3635.7144
$ws.Range("K27").Value = 1213.00002
$ws.Range("L27").Value = 10907.1432
$ws.Range("M27").Value = -1111.00002
$ws.Range("N27").Value = -11111.1432
# row 41 (Leve Item ID 4700)
$ws.Range("H41").Value = 1833.3334
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 3500
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 10500
$ws.Range("M41").Value = -2662
$ws.Range("N41").Value = -11176

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 2458.2
$ws.Range("J97").Value = 1189.5
$ws.Range("L97").Value = 1189.5
$ws.Range("N97").Value = -2181.5
# row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 1145.9445
$ws.Range("I122").Value = 911.53845
$ws.Range("K122").Value = 2734.61535
$ws.Range("M122").Value = -284.61535
# row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 66669804
$ws.Range("I132").Value = 90912680
$ws.Range("J132").Value = 1903.25
$ws.Range("K132").Value = 272738040
$ws.Range("L132").Value = 5709.75
$ws.Range("M132").Value = -272735510
$ws.Range("N132").Value = -10769.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 2 (Leve Item ID 2631)
$ws.Range("H2").Value = 15000
$ws.Range("J2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("N2").Value = -15224
# row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 1715.6666
$ws.Range("I68").Value = 1715.6666
$ws.Range("K68").Value = 1715.6666
$ws.Range("M68").Value = -966.6666
# row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 1715.6666
$ws.Range("I71").Value = 1715.6666
$ws.Range("K71").Value = 8578.333000000001
$ws.Range("M71").Value = -4834.333000000001
# row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 878.4545000000001
$ws.Range("J82").Value = 1371.75
$ws.Range("L82").Value = 1371.75
$ws.Range("N82").Value = -2093.75
# row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 878.4545000000001
$ws.Range("J85").Value = 1371.75
$ws.Range("L85").Value = 1371.75
$ws.Range("N85").Value = -3867.75
# row 101 (Leve Item ID 18549)
$ws.Range("H101").Value = 55723.668
$ws.Range("J101").Value = 55723.668
$ws.Range("L101").Value = 55723.668
$ws.Range("N101").Value = -62213.668
# row 108 (Leve Item ID 25655)
$ws.Range("H108").Value = 82000
$ws.Range("J108").Value = 82000
$ws.Range("L108").Value = 82000
$ws.Range("N108").Value = -89680
# row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 5196.273
$ws.Range("I122").Value = 5027.5884
$ws.Range("K122").Value = 15082.7652
$ws.Range("M122").Value = -12632.7652

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 13 (Leve Item ID 3008)
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
# row 103 (Leve Item ID 18548)
$ws.Range("H103").Value = 60142.6
$ws.Range("J103").Value = 60142.6
$ws.Range("L103").Value = 60142.6
$ws.Range("N103").Value = -62486.6
# row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 83333660
$ws.Range("J107").Value = 196.5
$ws.Range("L107").Value = 589.5
$ws.Range("N107").Value = -4429.5
# row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 5569.0356
$ws.Range("I126").Value = 5228.5
$ws.Range("K126").Value = 15685.5
$ws.Range("M126").Value = -13215.5
# row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1563.7391
$ws.Range("I132").Value = 1476.6487
$ws.Range("J132").Value = 1921.7778
$ws.Range("K132").Value = 4429.9461
$ws.Range("L132").Value = 5765.3334
$ws.Range("M132").Value = -1899.9461
$ws.Range("N132").Value = -10825.3334
# row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 821.2727
$ws.Range("I136").Value = 821.2727
$ws.Range("K136").Value = 2463.8181
$ws.Range("M136").Value = 86.18190000000004
